$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (shifts old I:K -> J:L), picking up the
# formatting of the neighbouring column automatically.
$ws.Columns("I:I").Insert()

# New header for the inserted column
$ws.Range("I1").Value = "ROI שכירות"

# New per-row values for the inserted column. These look like
# percentages, so Excel's auto-detection would otherwise coerce them to
# a numeric percent; force them back to plain text afterwards so they
# are stored as literal strings, same as the other text cells.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "0.00%"
$ws.Range("I2").ClearFormats()

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "0.00%"
$ws.Range("I3").ClearFormats()

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0.00%"
$ws.Range("I4").ClearFormats()

# Update the "summary" column (now K after the insert) text for rows 3-4
$ws.Range("K3").Value = "גבולי/לא משתלם"
$ws.Range("K4").Value = "גבולי/לא משתלם"
